# Insert a new column before the existing "x" column (column E) and
# give it the header "cs_relative_to", shifting the old "x" column to F.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Inserting at E shifts the current column E ("x" / axis / line) to F,
# and updates the sheet dimension from A1:E3 to A1:F3 automatically.
$ws.Range("E1").EntireColumn.Insert()

# New header cell.
$ws.Range("E1").Value = "cs_relative_to"

# The two data rows under the new header are blank text cells (matching
# the style of the other blank cells in the table, e.g. C2/D2/C3/D3).
# Assigning a single apostrophe produces an empty *text* cell (rather
# than an empty numeric/blank cell), and resetting the style afterwards
# clears the quote-prefix formatting so the cell ends up as a plain,
# unstyled, empty text cell.
$ws.Range("E2").Value = "'"
$ws.Range("E2").Style = "Normal"

$ws.Range("E3").Value = "'"
$ws.Range("E3").Style = "Normal"
